# The header row originally held text labels like "1961 [YR1961]" .. "2020
# [YR2020]" in F1:BM1 (E1 stays "1960 [YR1960]"). Replace them with actual
# numeric years 1960..2019 (shifted one column earlier than the old label
# series) so the row can drive a line graph, and left-align the new numeric
# year cells to match the old label look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$year = 1960
for ($col = 6; $col -le 65; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $year
    $cell.HorizontalAlignment = -4131  # xlLeft
    $year = $year + 1
}

# Restore the active selection to K5 like the edited workbook.
$ws.Range("K5").Select() | Out-Null
